$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New quarter of data (2022Q4) goes in row 61, directly below the existing
# last data row (60). Carry the same number formatting / styling as the
# row above by copying its formats down first.
$ws.Range("A60:BB60").Copy()
$ws.Range("A61:BB61").PasteSpecial(-4122)

$row = 61

$ws.Cells.Item($row, 1).Value = "2022Q4"
$ws.Cells.Item($row, 2).Value = 60148431
$ws.Cells.Item($row, 3).Value = 31804483
$ws.Cells.Item($row, 4).Value = 38607
$ws.Cells.Item($row, 5).Value = 1395053
$ws.Cells.Item($row, 6).Value = 1585924
$ws.Cells.Item($row, 7).Value = 55375
$ws.Cells.Item($row, 8).Value = 2882
$ws.Cells.Item($row, 9).Value = 139789
$ws.Cells.Item($row, 11).Value = 28586853
$ws.Cells.Item($row, 12).Value = 19058509
$ws.Cells.Item($row, 13).Value = 183129
$ws.Cells.Item($row, 14).Value = 113102
$ws.Cells.Item($row, 15).Value = 1232
$ws.Cells.Item($row, 16).Value = 2309
$ws.Cells.Item($row, 17).Value = 6183
$ws.Cells.Item($row, 18).Value = 22679
$ws.Cells.Item($row, 19).Value = 146553
$ws.Cells.Item($row, 20).Value = 116914
$ws.Cells.Item($row, 21).Value = 2236006
$ws.Cells.Item($row, 22).Value = 2455306
$ws.Cells.Item($row, 23).Value = 7755
$ws.Cells.Item($row, 24).Value = 16084
$ws.Cells.Item($row, 25).Value = 3898978
$ws.Cells.Item($row, 26).Value = 216132
$ws.Cells.Item($row, 27).Value = 74580
$ws.Cells.Item($row, 28).Value = 5363351
$ws.Cells.Item($row, 29).Value = 3026
$ws.Cells.Item($row, 30).Value = 683879
$ws.Cells.Item($row, 31).Value = 168417
$ws.Cells.Item($row, 32).Value = 30450
$ws.Cells.Item($row, 33).Value = 17494
$ws.Cells.Item($row, 34).Value = 3231
$ws.Cells.Item($row, 36).Value = 9381
$ws.Cells.Item($row, 37).Value = 4104
$ws.Cells.Item($row, 38).Value = 340003
$ws.Cells.Item($row, 39).Value = 464633
$ws.Cells.Item($row, 40).Value = 636985
$ws.Cells.Item($row, 41).Value = 36463
$ws.Cells.Item($row, 42).Value = 1800150
$ws.Cells.Item($row, 43).Value = 9100787
$ws.Cells.Item($row, 44).Value = 2428888
$ws.Cells.Item($row, 45).Value = 3266049
$ws.Cells.Item($row, 46).Value = 152917
$ws.Cells.Item($row, 47).Value = 475882
$ws.Cells.Item($row, 48).Value = 2040234
$ws.Cells.Item($row, 49).Value = 624187
$ws.Cells.Item($row, 50).Value = 100677
$ws.Cells.Item($row, 51).Value = 6838
$ws.Cells.Item($row, 52).Value = 5115
$ws.Cells.Item($row, 53).Value = 184652
$ws.Cells.Item($row, 54).Value = 184652

# Columns J (Trinidad and Tobago) and AI (Russia) have no data for this
# quarter yet, same as they are blank (formatted, but valueless) cells on
# other rows in this sheet.
$ws.Cells.Item($row, 10).ClearContents()
$ws.Cells.Item($row, 35).ClearContents()

$ws.Range("A61").Select()
